# Update crypto price/volume data as per the Dec 21 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text, matching the original inline-string cells.
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D29", "D30", "D31", "D33", "D34", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D50")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "43.704.69"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "2.236.84"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "269.98"
$ws.Range("E5").Value = "  +3.81%  "
$ws.Range("D6").Value = "92.92"
$ws.Range("E6").Value = "  +12.23%  "
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").Value = "0.622"
$ws.Range("E9").Value = "  +2.04%  "
$ws.Range("D10").Value = "46.92"
$ws.Range("E10").Value = "  +6.31%  "
$ws.Range("D11").Value = "0.0925"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "8.29"
$ws.Range("E12").Value = "  +17.47%  "
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "15.13"
$ws.Range("E14").Value = "  +3.42%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.564.85"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "2.238.87"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "0.802"
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("D18").Value = "43.641.20"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").Value = "0.0000104"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("D20").Value = "6.02"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").Value = "70.46"
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("D22").Value = "2.33"
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("D23").Value = "233.42"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "8.98"
$ws.Range("E24").Value = "  -3.69%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "11.36"
$ws.Range("E26").Value = "  +5.31%  "
$ws.Range("D27").Value = "2.50"
$ws.Range("E27").Value = "  +11.06%  "
$ws.Range("E28").Value = "  +5.35%  "
$ws.Range("D29").Value = "39.84"
$ws.Range("E29").Value = "  -3.98%  "
$ws.Range("D30").Value = "2.27"
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("D31").Value = "172.84"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  +4.10%  "
$ws.Range("D33").Value = "20.85"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("D34").Value = "5.49"
$ws.Range("E34").Value = "  +2.73%  "
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("E36").Value = "  -4.17%  "
$ws.Range("D37").Value = "0.0352"
$ws.Range("E37").Value = "  -3.12%  "
$ws.Range("D38").Value = "4.31"
$ws.Range("E38").Value = "  -4.05%  "
$ws.Range("D39").Value = "3.59"
$ws.Range("E39").Value = "  +21.47%  "
$ws.Range("D40").Value = "12.54"
$ws.Range("E40").Value = "  -6.66%  "
$ws.Range("D41").Value = "2.18"
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("D42").Value = "0.219"
$ws.Range("E42").Value = "  +8.69%  "
$ws.Range("D43").Value = "63.56"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "5.34"
$ws.Range("E44").Value = "  -3.94%  "
$ws.Range("D45").Value = "0.0990"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").Value = "100.62"
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("D48").Value = "1.16"
$ws.Range("E48").Value = "  +3.26%  "
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("D50").Value = "0.437"
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("D51").Value = "2.451.25"
$ws.Range("E51").Value = "  +0.22%  "

Write-Host "Updated $($ws.Name) with latest crypto data"
